$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 ("I0") and J1 ("IF"), styled like the other header cells (copy format from H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:J43
$iValues = @(6,9,6,4,7,5,7,6,9,9,6,8,9,5,7,7,6,3,3,6,9,4,7,9,7,9,6,7,9,9,2,9,5,6,8,8,7,5,6,6,7,3)
$jValues = @(7,9,6,5,8,5,8,6,9,9,6,9,9,5,8,8,8,5,4,6,9,5,7,9,7,9,6,7,9,9,4,9,5,7,8,8,7,6,7,6,7,3)

for ($k = 0; $k -lt $iValues.Length; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$k]
    $ws.Cells.Item($row, 10).Value = $jValues[$k]
}
